# Generate Report for Handoff
# Updates the GUID-based file names and timestamps across the Overview,
# zh-cn, and de-de sheets to reflect a newly generated handoff report.

$wb = $excel.ActiveWorkbook

$oldGuid = "2e85bbd5-b9ac-4191-be25-c098e0b84fc8"
$newGuid = "ec22fc48-0767-4738-93f9-0953479f8781"

$newZhXlf = "$newGuid.18cea1a9b852109002df29e9fac2fa31cd43af66.zh-cn.xlf"
$newDeXlf = "$newGuid.18cea1a9b852109002df29e9fac2fa31cd43af66.de-de.xlf"

# "Latest HO Xliff Generate Date" (Overview!G2) and the de-de sheet's
# "Latest Handoff Datetime" (de-de!H2) share the same shared-string value.
$newGenerateDate = "2016-09-07 13:48:22"
$newZhHandoffDate = "2016-09-07 13:47:58"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newGenerateDate
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate
foreach ($hl in $wsZh.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newGenerateDate
foreach ($hl in $wsDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
